$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Update the product name text (shared by both sheets' B1 cell) and the
# shortname value on the input sheet.
$ws1.Range("B1").Value = "2631-MS-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-VALIDATE-1st"
$ws2.Range("B1").Value = "2631-MS-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-VALIDATE-1st"
$ws1.Range("B2").Value = "263v"

# Remove the test-case inter-dependency: select B1 on the input sheet (instead
# of leaving the previous B18 selection) and make the output sheet the active
# tab instead of the input sheet.
$ws1.Range("B1").Select()
$ws2.Activate()
